# Fix "BREATING" -> "BREATHING" typo in the two headline occurrences.
#
#  1) Big centered title (sz 52): "BREATING CORRELATIONS PREVENTION"
#       -> "BREATHING CORRELATIONS PREVENTION"
#  2) Bold section heading (sz 24), right after the page break:
#       "BREATHING CORRELATIONS PREVENTION SECURITY SYSTEMS" already has the
#       correct spelling, but in the source XML "CORRELATIONS" and the
#       trailing space were split across two separate (but identically
#       formatted) runs; the edit simply re-merges them into a single run
#       containing "CORRELATIONS ".
#
# Track changes must be off while we edit so we get plain text edits
# instead of w:ins/w:del revision markup; restore the document's original
# setting afterwards so we don't leave an unrelated change behind in
# settings.xml.
$d = $word.ActiveDocument
$origTrackRevisions = $d.TrackRevisions
$d.TrackRevisions = $false

# --- 1) Fix the typo in the big title paragraph -----------------------
$fixedTitle = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range
    if ($rng.Text -like "*BREATING CORRELATIONS*") {
        $rng.Find.Execute("BREATING CORRELATIONS", $true, $false, $false,
                           $false, $false, $true, 1, $false,
                           "BREATHING CORRELATIONS", 2) | Out-Null
        $fixedTitle = $true
        break
    }
}

# --- 2) Re-merge "CORRELATIONS" + " " into one run in the heading -----
# Find the paragraph that reads "BREATHING CORRELATIONS PREVENTION SECURITY
# SYSTEMS" (the bold sz-24 heading under the page break) and touch only the
# part of it from "CORRELATIONS" onward, so the preceding "BREATHING " run
# (which carries the lastRenderedPageBreak marker) is left completely
# untouched.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range
    $text = $rng.Text
    if ($text -like "BREATHING CORRELATIONS PREVENTION SECURITY SYSTEMS*") {
        $relIdx = $text.IndexOf("CORRELATIONS")
        $subStart = $rng.Start + $relIdx
        $subRange = $d.Range($subStart, $rng.End)

        # Force a genuine (round-tripped) text change over just this
        # sub-range so the engine recombines "CORRELATIONS" and the
        # following space into a single run, without disturbing the
        # "BREATHING " run before it.
        $subRange.Find.Execute("CORRELATIONS ", $true, $false, $false,
                                $false, $false, $true, 1, $false,
                                "~~TMP~~CORRELATIONS~~TMP~~ ", 2) | Out-Null

        $subRange2 = $d.Range($subStart, $rng.End + 40)
        $subRange2.Find.Execute("~~TMP~~CORRELATIONS~~TMP~~ ", $true, $false,
                                 $false, $false, $false, $true, 1, $false,
                                 "CORRELATIONS ", 2) | Out-Null
        break
    }
}

# Restore the document's original TrackRevisions setting.
$d.TrackRevisions = $origTrackRevisions
